# Generate Report for Handback
# This script fills in the "handback" results for the file
# e7ce0a39-a3f3-4205-8085-0751ab805f5d.md (row 7) on both the
# zh-cn and de-de worksheets: the handback has now arrived, but it
# was generated from a stale source revision, so we record the
# target/handback file names, the handback datetime, a hyperlink to
# the (now outdated) target markdown file, and an error detail
# message describing the staleness.

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57b944171e6c82c73cd48002ffb2bacc537bc1f8/e2e/e7ce0a39-a3f3-4205-8085-0751ab805f5d.md"
$staleMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de38131bef4d4c0f00b864026bc1fef0efac467d/e2e/e7ce0a39-a3f3-4205-8085-0751ab805f5d.md"

$errorDetail = "The version of handback file is not the latest, current: " + $staleMdUrl + ", latest: " + $latestMdUrl + "."

# Hyperlink blue color used throughout this workbook (RGB 0x6495ED,
# expressed as BGR-packed OLE color for the Font.Color property).
$linkColor = 15570276

# ----- zh-cn sheet (row 7) -----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestMdUrl, "", "", "e7ce0a39-a3f3-4205-8085-0751ab805f5d.md")
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = $linkColor

$wsZh.Range("J7").Value = "e7ce0a39-a3f3-4205-8085-0751ab805f5d.b00579408c9064909c70d37ca47d7bf0b059fd93.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-07 09:25:48"
$wsZh.Range("P7").Value = $errorDetail

# ----- de-de sheet (row 7) -----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestMdUrl, "", "", "e7ce0a39-a3f3-4205-8085-0751ab805f5d.md")
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = $linkColor

$wsDe.Range("J7").Value = "e7ce0a39-a3f3-4205-8085-0751ab805f5d.b00579408c9064909c70d37ca47d7bf0b059fd93.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-07 09:25:56"
$wsDe.Range("P7").Value = $errorDetail
